$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("I2").Value = 0.7198826633506309
$ws.Range("J2").Value = 0.7198826633506309
$ws.Range("M2").Value = 50.86142466666666
$ws.Range("N2").Value = 152.584274
$ws.Range("O2").Value = 0.3434314568613803
$ws.Range("P2").Value = 0.3434314568613804
$ws.Range("Q2").Value = 3.062196841097778
$ws.Range("R2").Value = 27.55977156988
$ws.Range("S2").Value = 0.2472303518437578
$ws.Range("T2").Value = 0.2472303518437578

# Row 3
$ws.Range("I3").Value = 0.7198826633506309
$ws.Range("J3").Value = 0.7198826633506309
$ws.Range("M3").Value = 43.683024
$ws.Range("O3").Value = 0.294960761928139
$ws.Range("P3").Value = 0.294960761928139
$ws.Range("R3").Value = 23.67008338464
$ws.Range("S3").Value = 0.2123371388807601
$ws.Range("T3").Value = 0.2123371388807601

# Row 4
$ws.Range("I4").Value = 0.7198826633506309
$ws.Range("J4").Value = 0.7198826633506309
$ws.Range("M4").Value = 36.64360566666667
$ws.Range("N4").Value = 109.930817
$ws.Range("O4").Value = 0.2474285170192034
$ws.Range("P4").Value = 0.2474285170192035
$ws.Range("Q4").Value = 2.206189351837778
$ws.Range("R4").Value = 19.85570416654
$ws.Range("S4").Value = 0.1781194998206811
$ws.Range("T4").Value = 0.1781194998206811

# Row 5
$ws.Range("I5").Value = 0.7198826633506309
$ws.Range("J5").Value = 0.7198826633506309
$ws.Range("M5").Value = 16.90969166666667
$ws.Range("N5").Value = 50.729075
$ws.Range("O5").Value = 0.1141792641912772
$ws.Range("P5").Value = 0.1141792641912772
$ws.Range("Q5").Value = 1.018076169611111
$ws.Range("R5").Value = 9.162685526500001
$ws.Range("S5").Value = 0.08219567280543195
$ws.Range("T5").Value = 0.08219567280543197

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.02342733333333333
$ws.Range("H6").Value = 0.070282
$ws.Range("I6").Value = 0.2801173366493691
$ws.Range("J6").Value = 0.2801173366493691
$ws.Range("M6").Value = 50.86142466666666
$ws.Range("N6").Value = 152.584274
$ws.Range("O6").Value = 0.3434314568613803
$ws.Range("P6").Value = 0.3434314568613804
$ws.Range("Q6").Value = 1.191547549474222
$ws.Range("R6").Value = 10.723927945268
$ws.Range("S6").Value = 0.09620110501762254
$ws.Range("T6").Value = 0.09620110501762255

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.02342733333333333
$ws.Range("H7").Value = 0.070282
$ws.Range("I7").Value = 0.2801173366493691
$ws.Range("J7").Value = 0.2801173366493691
$ws.Range("M7").Value = 43.683024
$ws.Range("O7").Value = 0.294960761928139
$ws.Range("P7").Value = 0.294960761928139
$ws.Range("Q7").Value = 1.023376764256
$ws.Range("R7").Value = 9.210390878303999
$ws.Range("S7").Value = 0.0826236230473789
$ws.Range("T7").Value = 0.08262362304737891

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.02342733333333333
$ws.Range("H8").Value = 0.070282
$ws.Range("I8").Value = 0.2801173366493691
$ws.Range("J8").Value = 0.2801173366493691
$ws.Range("M8").Value = 36.64360566666667
$ws.Range("N8").Value = 109.930817
$ws.Range("O8").Value = 0.2474285170192034
$ws.Range("P8").Value = 0.2474285170192035
$ws.Range("Q8").Value = 0.8584619644882221
$ws.Range("R8").Value = 7.726157680394
$ws.Range("S8").Value = 0.06930901719852235
$ws.Range("T8").Value = 0.06930901719852235

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.02342733333333333
$ws.Range("H9").Value = 0.070282
$ws.Range("I9").Value = 0.2801173366493691
$ws.Range("J9").Value = 0.2801173366493691
$ws.Range("M9").Value = 16.90969166666667
$ws.Range("N9").Value = 50.729075
$ws.Range("O9").Value = 0.1141792641912772
$ws.Range("P9").Value = 0.1141792641912772
$ws.Range("Q9").Value = 0.3961489832388889
$ws.Range("R9").Value = 3.56534084915
$ws.Range("S9").Value = 0.03198359138584524
$ws.Range("T9").Value = 0.03198359138584525

